$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values
$ws.Range("C1").Value = "Score"
$ws.Range("D1").Value = "Salary"

# Add Score and Salary data for each row
$scores = @(73, 32, 46, 89, 4, 63, 25, 55, 95)
$salaries = @(100000, 35000, 25000, 55000, 15000, 60000, 30000, 150000, 80000)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $scores[$i]
    $ws.Cells.Item($row, 4).Value = $salaries[$i]
}

# Update selection to F9
$ws.Range("F9").Select() | Out-Null
